$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1879.8
$ws.Range("I40").Value = 1879.8
$ws.Range("K40").Value = 1879.8
$ws.Range("M40").Value = -1704.8
$ws.Range("H42").Value = 182.25
$ws.Range("I42").Value = 219.5
$ws.Range("J42").Value = 145
$ws.Range("K42").Value = 658.5
$ws.Range("L42").Value = 435
$ws.Range("M42").Value = -428.5
$ws.Range("N42").Value = -895
$ws.Range("H87").Value = 79999
$ws.Range("J87").Value = 79999
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495
$ws.Range("H90").Value = 79999
$ws.Range("J90").Value = 79999
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477
$ws.Range("H132").Value = 2342.7222
$ws.Range("I132").Value = 2332.182
$ws.Range("K132").Value = 6996.545999999999
$ws.Range("M132").Value = -4466.545999999999
$ws.Range("H138").Value = 7087.708
$ws.Range("I138").Value = 7670.737
$ws.Range("J138").Value = 6929.457
$ws.Range("K138").Value = 23012.211
$ws.Range("L138").Value = 20788.371
$ws.Range("M138").Value = -17872.211
$ws.Range("N138").Value = -31068.371

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 187.5
$ws.Range("I5").Value = 187.5
$ws.Range("K5").Value = 187.5
$ws.Range("M5").Value = -75.5
$ws.Range("H34").Value = 6999
$ws.Range("I34").Value = 6999
$ws.Range("K34").Value = 6999
$ws.Range("M34").Value = -6728
$ws.Range("H102").Value = 2051.6316
$ws.Range("J102").Value = 2166.5
$ws.Range("L102").Value = 2166.5
$ws.Range("N102").Value = -5410.5
$ws.Range("H110").Value = 10111.875
$ws.Range("I110").Value = 11082.5
$ws.Range("K110").Value = 11082.5
$ws.Range("M110").Value = -9037.5
$ws.Range("H122").Value = 1670501.4
$ws.Range("I122").Value = 2003601.6
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6010804.800000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6008354.800000001
$ws.Range("N122").Value = -19900
$ws.Range("H125").Value = 80000
$ws.Range("J125").Value = 80000
$ws.Range("L125").Value = 80000
$ws.Range("N125").Value = -89840
$ws.Range("H132").Value = 3450.8667
$ws.Range("I132").Value = 1810.8334
$ws.Range("J132").Value = 10011
$ws.Range("K132").Value = 5432.5002
$ws.Range("L132").Value = 30033
$ws.Range("M132").Value = -2902.5002
$ws.Range("N132").Value = -35093

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 187.5
$ws.Range("I4").Value = 187.5
$ws.Range("K4").Value = 187.5
$ws.Range("M4").Value = -72.5
$ws.Range("H61").Value = 48000
$ws.Range("J61").Value = 48000
$ws.Range("L61").Value = 48000
$ws.Range("N61").Value = -48626
$ws.Range("H86").Value = 4950.8
$ws.Range("I86").Value = 4542
$ws.Range("K86").Value = 4542
$ws.Range("M86").Value = -3419
$ws.Range("H89").Value = 4950.8
$ws.Range("I89").Value = 4542
$ws.Range("K89").Value = 22710
$ws.Range("M89").Value = -17094
$ws.Range("H94").Value = 1592.1
$ws.Range("I94").Value = 1703.1428
$ws.Range("J94").Value = 1333
$ws.Range("K94").Value = 1703.1428
$ws.Range("L94").Value = 1333
$ws.Range("M94").Value = -1252.1428
$ws.Range("N94").Value = -2235
$ws.Range("H105").Value = 3656.7715
$ws.Range("I105").Value = 2750.6924
$ws.Range("J105").Value = 6274.3335
$ws.Range("K105").Value = 2750.6924
$ws.Range("L105").Value = 6274.3335
$ws.Range("M105").Value = -1003.6924
$ws.Range("N105").Value = -9768.333500000001
$ws.Range("H107").Value = 4150.25
$ws.Range("I107").Value = 2798
$ws.Range("J107").Value = 9288.799999999999
$ws.Range("K107").Value = 2798
$ws.Range("L107").Value = 9288.799999999999
$ws.Range("M107").Value = -878
$ws.Range("N107").Value = -13128.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 819.6
$ws.Range("I16").Value = 699.6667
$ws.Range("K16").Value = 699.6667
$ws.Range("M16").Value = -412.6667
$ws.Range("H31").Value = 4289.2
$ws.Range("I31").Value = 3105.1365
$ws.Range("J31").Value = 5736.3887
$ws.Range("K31").Value = 3105.1365
$ws.Range("L31").Value = 5736.3887
$ws.Range("M31").Value = -2810.1365
$ws.Range("N31").Value = -6326.3887
$ws.Range("H34").Value = 4289.2
$ws.Range("I34").Value = 3105.1365
$ws.Range("J34").Value = 5736.3887
$ws.Range("K34").Value = 3105.1365
$ws.Range("L34").Value = 5736.3887
$ws.Range("M34").Value = -2903.1365
$ws.Range("N34").Value = -6140.3887
$ws.Range("H41").Value = 40000
$ws.Range("J41").Value = 40000
$ws.Range("L41").Value = 40000
$ws.Range("N41").Value = -40856
$ws.Range("H58").Value = 8569.571
$ws.Range("I58").Value = 7775
$ws.Range("J58").Value = 8702
$ws.Range("K58").Value = 7775
$ws.Range("L58").Value = 8702
$ws.Range("M58").Value = -7572
$ws.Range("N58").Value = -9108
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H107").Value = 749.4167
$ws.Range("I107").Value = 225.7
$ws.Range("K107").Value = 225.7
$ws.Range("M107").Value = 1694.3
$ws.Range("H113").Value = 819.6
$ws.Range("I113").Value = 699.6667
$ws.Range("K113").Value = 699.6667
$ws.Range("M113").Value = 1470.3333
$ws.Range("H132").Value = 2220.4814
$ws.Range("I132").Value = 2163.182
$ws.Range("K132").Value = 6489.545999999999
$ws.Range("M132").Value = -3959.545999999999
$ws.Range("H136").Value = 8569.571
$ws.Range("I136").Value = 7775
$ws.Range("J136").Value = 8702
$ws.Range("K136").Value = 23325
$ws.Range("L136").Value = 26106
$ws.Range("M136").Value = -20775
$ws.Range("N136").Value = -31206

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5125
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13372
$ws.Range("H65").Value = 5125
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 36000
$ws.Range("N65").Value = -42864
$ws.Range("H106").Value = 17394
$ws.Range("J106").Value = 17394
$ws.Range("L106").Value = 52182
$ws.Range("N106").Value = -54074
$ws.Range("H116").Value = 1349.5
$ws.Range("J116").Value = 1300
$ws.Range("L116").Value = 3900
$ws.Range("N116").Value = -10784
$ws.Range("H117").Value = 1724.6666
$ws.Range("J117").Value = 4197.5
$ws.Range("L117").Value = 12592.5
$ws.Range("N117").Value = -19476.5
$ws.Range("H129").Value = 1505.25
$ws.Range("J129").Value = 2030.5
$ws.Range("L129").Value = 6091.5
$ws.Range("N129").Value = -16091.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14000
$ws.Range("I15").Value = 14000
$ws.Range("K15").Value = 14000
$ws.Range("M15").Value = -13712
$ws.Range("H70").Value = 7000
$ws.Range("J70").Value = 7000
$ws.Range("L70").Value = 7000
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 7000
$ws.Range("J73").Value = 7000
$ws.Range("L73").Value = 7000
$ws.Range("N73").Value = -8872
$ws.Range("H81").Value = 14000
$ws.Range("I81").Value = 14000
$ws.Range("K81").Value = 14000
$ws.Range("M81").Value = -13002
$ws.Range("H84").Value = 14000
$ws.Range("I84").Value = 14000
$ws.Range("K84").Value = 42000
$ws.Range("M84").Value = -37008
$ws.Range("H102").Value = 2786.6667
$ws.Range("I102").Value = 2786.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2786.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1164.6667
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1579426.8
$ws.Range("I122").Value = 261998.25
$ws.Range("J122").Value = 3335998
$ws.Range("K122").Value = 785994.75
$ws.Range("L122").Value = 10007994
$ws.Range("M122").Value = -783544.75
$ws.Range("N122").Value = -10012894
$ws.Range("H126").Value = 4558
$ws.Range("I126").Value = 3851
$ws.Range("J126").Value = 4999.875
$ws.Range("K126").Value = 11553
$ws.Range("L126").Value = 14999.625
$ws.Range("M126").Value = -9083
$ws.Range("N126").Value = -19939.625
$ws.Range("H132").Value = 3099.1428
$ws.Range("I132").Value = 2127.3044
$ws.Range("K132").Value = 6381.9132
$ws.Range("M132").Value = -3851.9132

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H55").Value = 1307.421
$ws.Range("I55").Value = 974.3125
$ws.Range("K55").Value = 974.3125
$ws.Range("M55").Value = -801.3125
$ws.Range("H100").Value = 2867.9
$ws.Range("I100").Value = 2838
$ws.Range("J100").Value = 2987.5
$ws.Range("K100").Value = 2838
$ws.Range("L100").Value = 2987.5
$ws.Range("M100").Value = -2297
$ws.Range("N100").Value = -4069.5
$ws.Range("H132").Value = 5874.0835
$ws.Range("I132").Value = 3996.3333
$ws.Range("K132").Value = 11988.9999
$ws.Range("M132").Value = -9458.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2485.5715
$ws.Range("I107").Value = 399
$ws.Range("K107").Value = 1197
$ws.Range("M107").Value = 723
$ws.Range("H132").Value = 3367.4614
$ws.Range("I132").Value = 2797
$ws.Range("K132").Value = 8391
$ws.Range("M132").Value = -5861

